$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds numeric-looking IDs that must stay TEXT (shared string),
# matching the existing rows above. Apply a text format first so Excel
# stores "247" etc. as a string instead of auto-converting to a number,
# then drop back to the Normal style so no style index is left on the cell
# (matching rows 2-18, which have no s="" override on column A).
$ws.Range("A19:A36").NumberFormat = "@"

# Column G (Week) keeps the same text style used by the existing rows (s="1").
$ws.Range("G19:G36").NumberFormat = "@"

$ws.Range("A19").Value = "247"
$ws.Range("B19").Value = 1478.32
$ws.Range("C19").Value = 8.05798
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 184.76
$ws.Range("F19").Value = 13.3
$ws.Range("G19").Value = "Preinduction"

$ws.Range("A20").Value = "224"
$ws.Range("B20").Value = 2192.23
$ws.Range("C20").Value = 8.03838
$ws.Range("D20").Value = 25
$ws.Range("E20").Value = 275.08
$ws.Range("F20").Value = 3.82
$ws.Range("G20").Value = "Preinduction"

$ws.Range("A21").Value = "247"
$ws.Range("B21").Value = 2601.36
$ws.Range("C21").Value = 8.04627
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 325.58
$ws.Range("F21").Value = 16.68
$ws.Range("G21").Value = "1st week post induction"

$ws.Range("A22").Value = "224"
$ws.Range("B22").Value = 3181.81
$ws.Range("C22").Value = 8.99325
$ws.Range("D22").Value = 26
$ws.Range("E22").Value = 356.04
$ws.Range("F22").Value = 8.1
$ws.Range("G22").Value = "1st week post induction"

$ws.Range("A23").Value = "250"
$ws.Range("B23").Value = 3072.97
$ws.Range("C23").Value = 8.46036
$ws.Range("D23").Value = 25
$ws.Range("E23").Value = 365.64
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = "Preinduction"

$ws.Range("A24").Value = "256"
$ws.Range("B24").Value = 1473.44
$ws.Range("C24").Value = 7.68219
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = 193.68
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = "Preinduction"

$ws.Range("A25").Value = "262"
$ws.Range("B25").Value = 2108.71
$ws.Range("C25").Value = 6.27404
$ws.Range("D25").Value = 24
$ws.Range("E25").Value = 338.12
$ws.Range("F25").Value = 17.64
$ws.Range("G25").Value = "Preinduction"

$ws.Range("A26").Value = "271"
$ws.Range("B26").Value = 2971.73
$ws.Range("C26").Value = 9.35036
$ws.Range("D26").Value = 19
$ws.Range("E26").Value = 319.48
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "Preinduction"

$ws.Range("A27").Value = "281"
$ws.Range("B27").Value = 667.434
$ws.Range("C27").Value = 5.36436
$ws.Range("D27").Value = 29
$ws.Range("E27").Value = 126.8
$ws.Range("F27").Value = 32.92
$ws.Range("G27").Value = "Preinduction"

$ws.Range("A28").Value = "250"
$ws.Range("B28").Value = 1887.67
$ws.Range("C28").Value = 8.32011
$ws.Range("D28").Value = 15
$ws.Range("E28").Value = 228.16
$ws.Range("F28").Value = 18.46
$ws.Range("G28").Value = "1st week post induction"

$ws.Range("A29").Value = "256"
$ws.Range("B29").Value = 1196.44
$ws.Range("C29").Value = 8.51315
$ws.Range("D29").Value = 17
$ws.Range("E29").Value = 142.1
$ws.Range("F29").Value = 10.88
$ws.Range("G29").Value = "1st week post induction"

$ws.Range("A30").Value = "262"
$ws.Range("B30").Value = 2217.97
$ws.Range("C30").Value = 7.38929
$ws.Range("D30").Value = 29
$ws.Range("E30").Value = 302.78
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = "1st week post induction"

$ws.Range("A31").Value = "271"
$ws.Range("B31").Value = 2299.56
$ws.Range("C31").Value = 6.56005
$ws.Range("D31").Value = 17
$ws.Range("E31").Value = 352.04
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = "1st week post induction"

$ws.Range("A32").Value = "281"
$ws.Range("B32").Value = 836.462
$ws.Range("C32").Value = 5.485
$ws.Range("D32").Value = 24
$ws.Range("E32").Value = 154.52
$ws.Range("F32").Value = 21.92
$ws.Range("G32").Value = "1st week post induction"

$ws.Range("A33").Value = "250"
$ws.Range("B33").Value = 4078.39
$ws.Range("C33").Value = 11.3062
$ws.Range("D33").Value = 20
$ws.Range("E33").Value = 362.5
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = "2nd week post induction"

$ws.Range("A34").Value = "256"
$ws.Range("B34").Value = 2413.92
$ws.Range("C34").Value = 9.23601
$ws.Range("D34").Value = 20
$ws.Range("E34").Value = 263.08
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = "2nd week post induction"

$ws.Range("A35").Value = "262"
$ws.Range("B35").Value = 1531.33
$ws.Range("C35").Value = 6.90908
$ws.Range("D35").Value = 33
$ws.Range("E35").Value = 224.52
$ws.Range("F35").Value = 11.22
$ws.Range("G35").Value = "2nd week post induction"

$ws.Range("A36").Value = "281"
$ws.Range("B36").Value = 1893.1
$ws.Range("C36").Value = 6.62201
$ws.Range("D36").Value = 35
$ws.Range("E36").Value = 289.02
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = "2nd week post induction"

# Restore column A to the default (unstyled) cell style; column G intentionally
# keeps the "@" text style applied above, matching the target workbook.
$ws.Range("A19:A36").Style = "Normal"

# Mirror the author's final navigation/selection state: scrolled down so
# row 8 is the top visible row, with the active selection on H37.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H37").Select()
